$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, copying the format of the existing header cells (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Rows where the new "Save" value is 1 (everything else is 0)
$saveOnes = @(11, 21, 27)

for ($r = 2; $r -le 32; $r++) {
    if ($saveOnes -contains $r) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
